# Update countries & provincias Spain
# Refreshes the "Pais" covid-tracker sheet: bumps the "Datos actualizados" timestamp,
# updates the daily statistics for a handful of countries, and reflects the resulting
# re-sort (two pairs of countries swap ranking positions: Ucrania/Israel and
# Kirguistan/Ghana, taking their respective stat rows along with them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 3 de Agosto de 2020 a las 08:36"

# --- Row 31: Ecuador (stats-only update) -------------------------------
$ws.Cells.Item(31, 2).Value = 86524
$ws.Cells.Item(31, 5).Value = 44561
$ws.Cells.Item(31, 8).Value = 5750

# --- Rows 36-38: Ucrania overtakes Israel & Republica Dominicana -------
# New ranking order: Ucrania, Israel, Republica Dominicana
$ws.Cells.Item(36, 1).Value = "Ucrania"
$ws.Cells.Item(36, 2).Value = 73158
$ws.Cells.Item(36, 3).Value = 990
$ws.Cells.Item(36, 4).Value = 39876
$ws.Cells.Item(36, 5).Value = 31544
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 13
$ws.Cells.Item(36, 8).Value = 1738

$ws.Cells.Item(37, 1).Value = "Israel"
$ws.Cells.Item(37, 2).Value = 73025
$ws.Cells.Item(37, 3).Value = 210
$ws.Cells.Item(37, 4).Value = 47511
$ws.Cells.Item(37, 5).Value = 24978
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 8).Value = 536

$ws.Cells.Item(38, 1).Value = "Republica Dominicana"
$ws.Cells.Item(38, 2).Value = 72243
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(38, 4).Value = 38244
$ws.Cells.Item(38, 5).Value = 32821
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(38, 8).Value = 1178

# --- Row 53: Armenia (stats-only update) -------------------------------
$ws.Cells.Item(53, 2).Value = 39102
$ws.Cells.Item(53, 3).Value = 52
$ws.Cells.Item(53, 4).Value = 29861
$ws.Cells.Item(53, 5).Value = 8479
$ws.Cells.Item(53, 7).Value = 8
$ws.Cells.Item(53, 8).Value = 762

# --- Rows 54-55: Kirguistan overtakes Ghana ----------------------------
$ws.Cells.Item(54, 1).Value = "Kirguistan"
$ws.Cells.Item(54, 2).Value = 37129
$ws.Cells.Item(54, 3).Value = 410
$ws.Cells.Item(54, 4).Value = 27927
$ws.Cells.Item(54, 5).Value = 7782
$ws.Cells.Item(54, 7).Value = 11
$ws.Cells.Item(54, 8).Value = 1420

$ws.Cells.Item(55, 1).Value = "Ghana"
$ws.Cells.Item(55, 2).Value = 37014
$ws.Cells.Item(55, 4).Value = 33365
$ws.Cells.Item(55, 5).Value = 3467
$ws.Cells.Item(55, 8).Value = 182

# --- Row 63: Uzbekistan (stats-only update) -----------------------------
$ws.Cells.Item(63, 2).Value = 25828
$ws.Cells.Item(63, 3).Value = 492
$ws.Cells.Item(63, 5).Value = 9166
$ws.Cells.Item(63, 7).Value = 4
$ws.Cells.Item(63, 8).Value = 155

# --- Row 72: Australia (stats-only update) ------------------------------
$ws.Cells.Item(72, 2).Value = 18318
$ws.Cells.Item(72, 3).Value = 395
$ws.Cells.Item(72, 4).Value = 10622
$ws.Cells.Item(72, 5).Value = 7475
$ws.Cells.Item(72, 7).Value = 13
$ws.Cells.Item(72, 8).Value = 221

# --- Row 73: El Salvador (stats-only update) -----------------------------
$ws.Cells.Item(73, 4).Value = 8649
$ws.Cells.Item(73, 5).Value = 8322
$ws.Cells.Item(73, 7).Value = 10
$ws.Cells.Item(73, 8).Value = 477

# --- Row 106: Hungria (stats-only update) --------------------------------
$ws.Cells.Item(106, 2).Value = 4544
$ws.Cells.Item(106, 3).Value = 9
$ws.Cells.Item(106, 4).Value = 3413
$ws.Cells.Item(106, 5).Value = 534

# --- Row 145: Georgia (stats-only update) --------------------------------
$ws.Cells.Item(145, 2).Value = 1179
$ws.Cells.Item(145, 3).Value = 2
$ws.Cells.Item(145, 4).Value = 959
$ws.Cells.Item(145, 5).Value = 203

# --- Row 176: Camboya (stats-only update) --------------------------------
$ws.Cells.Item(176, 4).Value = 197
$ws.Cells.Item(176, 5).Value = 43
